$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("A11").Value = "HMS Heart Prediction"
$ws.Range("B11").Value = "POST"
$ws.Range("C11").Value = "http://127.0.0.1:8000/user/predictHeart/"
$ws.Range("D11").Value = '{"age":45,"sex":1,"cp":1,"trestbps":145,"chol":220,"fbs":1,"restecg":2,"thalach":150,"exang":1, "oldpeak":2.34,"slope":3,"ca":3.2, "thal":5.0}'
$ws.Range("E11").Value = '{"result": "The person does not have a Heart Disease", "tips": ["Maintain a balanced diet with fruits, vegetables, lean
proteins, and healthy fats.", "Stay physically active to maintain a healthy weight and cardiovascular fitness.", "Avoid
smoking and seek help to quit if needed.", "Limit alcohol consumption to moderate levels.", "Practice stress-relief
techniques to reduce the impact of stress.", "Schedule regular health check-ups for early detection and prevention."],
"youtube_links": {"heart_disease_management": "https://www.youtube.com/watch?v=IMBpwpf5crU", "heart_disease_prevention":
"https://www.youtube.com/watch?v=B6UYNZLpAMs"}}'

$ws.Hyperlinks.Add($ws.Range("C11"), "http://127.0.0.1:8000/user/predictHeart/") | Out-Null

$ws.Range("A11:E11").HorizontalAlignment = -4108
$ws.Range("A11:E11").VerticalAlignment = -4108
$ws.Range("A11:E11").WrapText = $true
$ws.Range("A11").Font.Bold = $true
$ws.Range("B11").Font.Bold = $true

$ws.Rows.Item(11).RowHeight = 244.8

$excel.ActiveWindow.ScrollRow = 11
$ws.Range("C15").Select()
